# Applies the "Add new columns of data to existing worksheet" edit.
#
# Target state (per the commit diff):
#  - Sheet2        (tab 1, rId1) -> rewritten to a 9-col x 6-row table (A1:I6)
#  - Sheet2 (2)    (tab 2, rId2) -> renamed to "Sheet 2 copy" and rewritten to
#                                    a 6-col x 6-row table (A1:F6)
#  - Sheet1        (tab 3, rId3) -> left untouched

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1 ("Sheet2") - new 9 column layout
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet2")

# Grab a cell that already carries the bold/border/center style (style index 1)
# before we wipe the sheet, so we can clone that exact formatting later via
# PasteSpecial instead of accidentally creating brand new font/xf entries.
$styleSource = $wb.Worksheets.Item("Sheet1").Range("A1")

$ws1.Cells.Clear()

$ws1.Range("A1").Value = " Old Column"
$ws1.Range("B1").Value = "New Column D"
$ws1.Range("C1").Value = "New Column C"
$ws1.Range("D1").Value = "New Column B"
$ws1.Range("E1").Value = "New Column A"
$ws1.Range("F1").Value = "Old 2"
$ws1.Range("G1").Value = "Old 3"
$ws1.Range("H1").Value = "Old 4"
$ws1.Range("I1").Value = "Old 5"

$ws1.Range("A2").Value = 1
$ws1.Range("B2").Value = 2
$ws1.Range("C2").Value = 2
$ws1.Range("D2").Value = 1
$ws1.Range("E2").Value = 1
$ws1.Range("F2").Value = 1
$ws1.Range("G2").Value = 10
$ws1.Range("H2").Value = 10
$ws1.Range("I2").Value = 10

$ws1.Range("A3").Value = 2
$ws1.Range("B3").Value = 3
$ws1.Range("C3").Value = 3
$ws1.Range("D3").Value = 2
$ws1.Range("E3").Value = 2
$ws1.Range("F3").Value = 2
$ws1.Range("G3").Value = 20
$ws1.Range("H3").Value = 20
$ws1.Range("I3").Value = 20

$ws1.Range("A4").Value = 3
$ws1.Range("B4").Value = 4
$ws1.Range("C4").Value = 4
$ws1.Range("D4").Value = 3
$ws1.Range("E4").Value = 3
$ws1.Range("F4").Value = 3
$ws1.Range("G4").Value = 30
$ws1.Range("H4").Value = 30
$ws1.Range("I4").Value = 30

$ws1.Range("A5").Value = 4
$ws1.Range("B5").Value = 5
$ws1.Range("C5").Value = 5
$ws1.Range("D5").Value = 4
$ws1.Range("E5").Value = 4
$ws1.Range("F5").Value = 4

$ws1.Range("A6").Value = 5
$ws1.Range("B6").Value = 6
$ws1.Range("C6").Value = 6
$ws1.Range("D6").Value = 5
$ws1.Range("E6").Value = 5
$ws1.Range("F6").Value = 5

# Re-apply the custom column widths that were already present on columns A-C,
# then extend the same sizing scheme to the new D-F columns. (The ColumnWidth
# COM property is quantized to ~1/6-character steps by the host, same as real
# Excel's MDW-based rounding, so these land on the closest representable
# width to the original 19.5 / 18.1640625 / 21.5 / 30.6640625 / 18.83203125.)
$ws1.Range("A1").ColumnWidth = 18.6
$ws1.Range("B1").ColumnWidth = 17.3
$ws1.Range("C1").ColumnWidth = 20.6
$ws1.Range("D1").ColumnWidth = 29.8
$ws1.Range("E1").ColumnWidth = 17.95
$ws1.Range("F1").ColumnWidth = 20.6

# Re-apply the bold/border/center header formatting to G1:I1 (the old "Old
# 3/4/5" header cells), matching style index 1 in the original stylesheet.
$styleSource.Copy()
$ws1.Range("G1:I1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Range("B10").Select()

# ---------------------------------------------------------------------------
# Sheet 2 ("Sheet2 (2)" -> "Sheet 2 copy") - new 6 column layout
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2 (2)")
$ws2.Name = "Sheet 2 copy"

$ws2.Cells.Clear()

$ws2.Range("A1").Value = "Old 1"
$ws2.Range("B1").Value = "New Column A"
$ws2.Range("C1").Value = "Old 2"
$ws2.Range("D1").Value = "Old 3"
$ws2.Range("E1").Value = "Old 4"
$ws2.Range("F1").Value = "Old 5"

$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = 1
$ws2.Range("D2").Value = 10
$ws2.Range("E2").Value = 10
$ws2.Range("F2").Value = 10

$ws2.Range("A3").Value = 2
$ws2.Range("B3").Value = 2
$ws2.Range("C3").Value = 2
$ws2.Range("D3").Value = 20
$ws2.Range("E3").Value = 20
$ws2.Range("F3").Value = 20

$ws2.Range("A4").Value = 3
$ws2.Range("B4").Value = 3
$ws2.Range("C4").Value = 3
$ws2.Range("D4").Value = 30
$ws2.Range("E4").Value = 30
$ws2.Range("F4").Value = 30

$ws2.Range("A5").Value = 4
$ws2.Range("B5").Value = 4
$ws2.Range("C5").Value = 4

$ws2.Range("A6").Value = 5
$ws2.Range("B6").Value = 5
$ws2.Range("C6").Value = 5

# Restore the original A-C column widths (Clear() wipes per-column sizing).
$ws2.Range("A1").ColumnWidth = 18.666666666666668
$ws2.Range("B1").ColumnWidth = 17.330729166666668
$ws2.Range("C1").ColumnWidth = 20.666666666666668

# Re-apply the bold/border/center header formatting to D1:F1 (the old "Old
# 3/4/5" header cells), matching style index 1 in the original stylesheet.
$styleSource.Copy()
$ws2.Range("D1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws2.Range("F37").Select()
